#
# 10 29 2016 yyy
#
# Re-creates the author's edit:
#   - workbook.xml: drop the saved firstSheet/activeTab window state
#   - pay_order sheet: activate it, change the selection, widen/extend
#     columns P..W, and fill in previously-empty cells in rows 2,3,5
#   - pay_settle_result sheet: no longer the active (tabSelected) sheet
#

$wb = $excel.ActiveWorkbook

$payOrder  = $wb.Worksheets.Item("pay_order")
$payResult = $wb.Worksheets.Item("pay_settle_result")

# --- pay_order: become the active sheet, select W2:W4 -----------------
$payOrder.Select()
$payOrder.Range("W2:W4").Select()

# --- pay_order: column widths -----------------------------------------
# Column P (16) widens and loses its "best fit" flag.
$payOrder.Columns.Item(16).ColumnWidth = 26.285714285714285
# New columns Q..W (17-23) get explicit widths.
$payOrder.Columns.Item(17).ColumnWidth = 21.428571428571427
$payOrder.Columns.Item(18).ColumnWidth = 19.428571428571427
$payOrder.Columns.Item(19).ColumnWidth = 20.428571428571427
$payOrder.Columns.Item(20).ColumnWidth = 24.571428571428573
$payOrder.Columns.Item(21).ColumnWidth = 14.285714285714286
$payOrder.Columns.Item(22).ColumnWidth = 13.0
$payOrder.Columns.Item(23).ColumnWidth = 13.857142857142858

# --- pay_order: row 2 & 3, fill in S (hyperlink-styled) and T ---------
$payOrder.Range("S2").Copy()
$payOrder.Range("S3").PasteSpecial(-4122)
$payOrder.Range("S2").Value = "http://www.qq.com"
$payOrder.Range("S3").Value = "http://www.qq.com"
$payOrder.Range("T2").Value = "http://www.baidu.com"
$payOrder.Range("T3").Value = "http://www.baidu.com"

# --- pay_order: row 5, P/Q/R switch from raw dates to text, plus S5/T5
$payOrder.Range("P5:R5").NumberFormat = "@"
$payOrder.Range("P5").Value = "2016-07-18 11:47:30"
$payOrder.Range("Q5").Value = "2016-07-18 11:46:16"
$payOrder.Range("R5").Value = "2016-07-18 11:47:31"

$payOrder.Range("S2").Copy()
$payOrder.Range("S5").PasteSpecial(-4122)
$payOrder.Range("S5").Value = "http://www.qq.com"
$payOrder.Range("T5").Value = "http://www.baidu.com"

# --- pay_settle_result: no longer tab-selected -------------------------
# Selecting pay_order above already moves tabSelected off this sheet;
# its own selection (A7) is left untouched.
